# Update cryptos list - GitHub Actions scheduled data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
}

# Row 2 - Bitcoin
Set-Text "D2" "69.831.77"
Set-Text "E2" "  +2.62%  "

# Row 3 - Ethereum
Set-Text "D3" "3.947.51"
Set-Text "E3" "  +1.16%  "

# Row 4 - TetherUSD
Set-Text "E4" "  +0.10%  "

# Row 5 - BNB
Set-Text "D5" "527.72"
Set-Text "E5" "  +8.28%  "

# Row 6 - Solana
Set-Text "D6" "146.78"
Set-Text "E6" "  +0.42%  "

# Row 7 - XRP
Set-Text "E7" "  +0.45%  "

# Row 8 - USDC
Set-Text "D8" "0.998"
Set-Text "E8" "  +0.05%  "

# Row 9 - Cardano
Set-Text "E9" "  +0.83%  "

# Row 10 - Dogecoin
Set-Text "E10" "  +5.54%  "

# Row 11 - ShibaInu
Set-Text "E11" "  -0.20%  "

# Row 12 - Avalanche
Set-Text "D12" "42.90"
Set-Text "E12" "  +0.38%  "

# Row 13 - Polkadot
Set-Text "D13" "10.52"
Set-Text "E13" "  -1.41%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-Text "D14" "4.585.63"
Set-Text "E14" "  +1.66%  "

# Row 15 - WrappedEther
Set-Text "D15" "3.946.38"
Set-Text "E15" "  +0.68%  "

# Row 16 - Uniswap
Set-Text "D16" "14.17"
Set-Text "E16" "  +1.04%  "

# Row 17 - TRON
Set-Text "E17" "  -0.07%  "

# Row 18 - Polygon
Set-Text "E18" "  +7.46%  "

# Row 19 - Chainlink
Set-Text "D19" "19.91"
Set-Text "E19" "  +0.84%  "

# Row 20 - WrappedBTC
Set-Text "D20" "69.769.92"
Set-Text "E20" "  +2.34%  "

# Row 21 - BitcoinCash
Set-Text "D21" "436.22"
Set-Text "E21" "  +1.75%  "

# Row 22 - ImmutableX
Set-Text "D22" "3.41"
Set-Text "E22" "  -3.80%  "

# Row 23 - InternetComputer(DFINITY)
Set-Text "D23" "14.57"
Set-Text "E23" "  -2.47%  "

# Row 24 - Litecoin
Set-Text "D24" "88.60"
Set-Text "E24" "  +1.64%  "

# Row 25 - was PancakeSwap, now RenderToken
Set-Text "B25" "RenderToken"
Set-Text "C25" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-Text "D25" "11.96"
Set-Text "E25" "  +6.27%  "

# Row 26 - was RenderToken, now PancakeSwap
Set-Text "B26" "PancakeSwap"
Set-Text "C26" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-Text "D26" "4.03"
Set-Text "E26" "  +11.53%  "

# Row 27 - Filecoin
Set-Text "E27" "  -1.93%  "

# Row 28 - EthereumClassic
Set-Text "D28" "36.81"
Set-Text "E28" "  -3.37%  "

# Row 29 - LEO
Set-Text "D29" "5.67"
Set-Text "E29" "  -1.13%  "

# Row 30 - Bittensor
Set-Text "D30" "703.22"
Set-Text "E30" "  -2.34%  "

# Row 31 - Cosmos
Set-Text "D31" "13.37"
Set-Text "E31" "  -2.36%  "

# Row 32 - Hedera
Set-Text "E32" "  -2.00%  "

# Row 33 - Toncoin
Set-Text "E33" "  -1.04%  "

# Row 34 - OKB
Set-Text "D34" "67.00"
Set-Text "E34" "  +11.18%  "

# Row 35 - TheGraph
Set-Text "D35" "0.442"
Set-Text "E35" "  +8.50%  "

# Row 36 - PEPE
$subscript3 = [string][char]0x2083
Set-Text "D36" "0.0${subscript3}0876"
Set-Text "E36" "  +1.34%  "

# Row 37 - NEARProtocol
Set-Text "D37" "6.03"
Set-Text "E37" "  -4.22%  "

# Row 38 - InjectiveProtocol
Set-Text "D38" "40.44"
Set-Text "E38" "  -2.83%  "

# Row 39 - Kaspa
Set-Text "E39" "  +0.33%  "

# Row 40 - Dai
Set-Text "D40" "0.998"
Set-Text "E40" "  -0.11%  "

# Row 41 - FirstDigitalUSD
Set-Text "E41" "  +0.12%  "

# Row 42 - VeChain
Set-Text "E42" "  +1.50%  "

# Row 43 - Fetch.AI
Set-Text "D43" "2.86"
Set-Text "E43" "  -3.32%  "

# Row 44 - WEMIXToken
Set-Text "D44" "3.11"
Set-Text "E44" "  +7.39%  "

# Row 45 - ThetaToken
Set-Text "D45" "3.01"
Set-Text "E45" "  -3.85%  "

# Row 46 - was ApeXProtocol, now Stacks
Set-Text "B46" "Stacks"
Set-Text "C46" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-Text "D46" "3.18"
Set-Text "E46" "  +13.22%  "

# Row 47 - was Stellar, now ApeXProtocol
Set-Text "B47" "ApeXProtocol"
Set-Text "C47" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-Text "D47" "3.41"
Set-Text "E47" "  +2.99%  "

# Row 48 - was Stacks, now Stellar
Set-Text "B48" "Stellar"
Set-Text "C48" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-Text "D48" "0.143"
Set-Text "E48" "  +1.65%  "

# Row 49 - BabyDogeCoin
Set-Text "E49" "  +8.92%  "

# Row 50 - LidoDAOToken
Set-Text "D50" "3.36"
Set-Text "E50" "  -0.96%  "

# Row 51 - ARBITRUM
Set-Text "D51" "2.11"
Set-Text "E51" "  -0.55%  "
